$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D21: remove the "adverse treatment outcome [PATO:0000011]; " prefix
$ws.Range("D21").Value = "physical quality [PATO:0001018]; concentration [PATO:0000033]"

# Add new row 28 with the imported NCIT ontology entry
$ws.Range("A28").Value = "NCIT"
$ws.Range("B28").Value = "http://purl.obolibrary.org/obo/ncit.owl"
$ws.Range("C28").Value = "entity [BFO:0000001]"
$ws.Range("D28").Value = "agitation [NCIT:C79530]; thalamus [NCIT:C12459]"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "all"
$ws.Range("F28").Value = "'"
$ws.Range("F28").Style = "Normal"
